$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the "Tropopause " header text (remove trailing space)
$ws.Range("E1").Value = "Tropopause"

# Fill column E (Tropopause height) with 12000 for all data rows (3-52)
for ($r = 3; $r -le 52; $r++) {
    $ws.Cells.Item($r, 5).Value = 12000
}
